$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.536181092262268
$ws.Range("B1").Value = 2.13808012008667
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.378749966621399
$ws.Range("E1").Value = 0.6528251171112061
